$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for the two records (currently on row 3 and row 6) need to
# be swapped: row 3 should take on the values currently in row 6, and row 6
# should take on the values currently in row 3. Only columns A, B, E, F, G,
# H, Q, R and AC differ between the two records; the rest are identical so
# swapping them is a no-op, but we swap the full set of differing columns
# here (including the AC "Publik kommentar" text, which only exists on one
# of the two rows).

# Capture current (pre-swap) values from row 3
# (use Value2, since Value can return a COM property wrapper in this
# environment instead of the actual scalar data)
$A3 = $ws.Range("A3").Value2
$B3 = $ws.Range("B3").Value2
$E3 = $ws.Range("E3").Value2
$F3 = $ws.Range("F3").Value2
$G3 = $ws.Range("G3").Value2
$H3 = $ws.Range("H3").Value2
$Q3 = $ws.Range("Q3").Value2
$R3 = $ws.Range("R3").Value2
$AC3 = $ws.Range("AC3").Value2

# Capture current (pre-swap) values from row 6
$A6 = $ws.Range("A6").Value2
$B6 = $ws.Range("B6").Value2
$E6 = $ws.Range("E6").Value2
$F6 = $ws.Range("F6").Value2
$G6 = $ws.Range("G6").Value2
$H6 = $ws.Range("H6").Value2
$Q6 = $ws.Range("Q6").Value2
$R6 = $ws.Range("R6").Value2
$AC6 = $ws.Range("AC6").Value2

# Write row 6's original values into row 3
$ws.Range("A3").Value2 = $A6
$ws.Range("B3").Value2 = $B6
$ws.Range("E3").Value2 = $E6
$ws.Range("F3").Value2 = $F6
$ws.Range("G3").Value2 = $G6
$ws.Range("H3").Value2 = $H6
$ws.Range("Q3").Value2 = $Q6
$ws.Range("R3").Value2 = $R6
$ws.Range("AC3").ClearContents()

# Write row 3's original values into row 6
$ws.Range("A6").Value2 = $A3
$ws.Range("B6").Value2 = $B3
$ws.Range("E6").Value2 = $E3
$ws.Range("F6").Value2 = $F3
$ws.Range("G6").Value2 = $G3
$ws.Range("H6").Value2 = $H3
$ws.Range("Q6").Value2 = $Q3
$ws.Range("R6").Value2 = $R3
$ws.Range("AC6").Value2 = $AC3
